$wb = $excel.ActiveWorkbook

# --- Sheets ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Color used for the existing "HyperLink" cell style (rgb FF6495ED == RGB(100,149,237))
$hyperlinkColor = 15570276

# --- 1. Overview sheet: update the per-locale handback status text ---
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the zh-cn / de-de status columns on the Overview sheet to fit the new text
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- 2. zh-cn sheet ---
# Widen the "Status" column (C) and the "Latest Handback File" column (J)
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# Latest Target File (I) / Latest Handback File (J) for both data rows
$wsZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# Rebuild the hyperlinks so the relationship ids come out in the same order as Excel would
# produce them (A2, I2, A3, I3) -- this also gives I2 / I3 their "a.md" hyperlink.
$wsZhCn.Hyperlinks.Delete()

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec3d012e024eb1b9a15f8078f92c027aae818f63/e2e/a.md", "", "", "a.md")
$wsZhCn.Range("A2").Font.Underline = $true
$wsZhCn.Range("A2").Font.Color = $hyperlinkColor

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec3d012e024eb1b9a15f8078f92c027aae818f63/e2e/a.md", "", "", "a.md")
$wsZhCn.Range("I2").Font.Underline = $true
$wsZhCn.Range("I2").Font.Color = $hyperlinkColor

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec3d012e024eb1b9a15f8078f92c027aae818f63/e2e/b.md", "", "", "b.md")
$wsZhCn.Range("A3").Font.Underline = $true
$wsZhCn.Range("A3").Font.Color = $hyperlinkColor

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec3d012e024eb1b9a15f8078f92c027aae818f63/e2e/a.md", "", "", "a.md")
$wsZhCn.Range("I3").Font.Underline = $true
$wsZhCn.Range("I3").Font.Color = $hyperlinkColor

# Handback datetime unchanged in this run (still "0001-01-01 00:00:00" conceptually, but the
# underlying shared string value moves forward to the latest handoff run time)
$wsZhCn.Range("K2").Value = "2016-08-16 18:32:20"
$wsZhCn.Range("K3").Value = "2016-08-16 18:32:20"

# --- 3. de-de sheet ---
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDeDe.Hyperlinks.Delete()

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec3d012e024eb1b9a15f8078f92c027aae818f63/e2e/a.md", "", "", "a.md")
$wsDeDe.Range("A2").Font.Underline = $true
$wsDeDe.Range("A2").Font.Color = $hyperlinkColor

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec3d012e024eb1b9a15f8078f92c027aae818f63/e2e/a.md", "", "", "a.md")
$wsDeDe.Range("I2").Font.Underline = $true
$wsDeDe.Range("I2").Font.Color = $hyperlinkColor

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec3d012e024eb1b9a15f8078f92c027aae818f63/e2e/b.md", "", "", "b.md")
$wsDeDe.Range("A3").Font.Underline = $true
$wsDeDe.Range("A3").Font.Color = $hyperlinkColor

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec3d012e024eb1b9a15f8078f92c027aae818f63/e2e/a.md", "", "", "a.md")
$wsDeDe.Range("I3").Font.Underline = $true
$wsDeDe.Range("I3").Font.Color = $hyperlinkColor

# de-de handback has its own, later, completion timestamp
$wsDeDe.Range("K2").Value = "2016-08-16 18:32:28"
$wsDeDe.Range("K3").Value = "2016-08-16 18:32:28"
